$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from AC1 (format only) to new header cells AD1:AF1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set new header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for each data row (2-43)
for ($row = 2; $row -le 43; $row++) {
    $ws.Range("AD" + $row).Value = 90
    $ws.Range("AE" + $row).Value = 72
    $ws.Range("AF" + $row).Value = 0
}

Write-Output "done"
